# "names are removed from the results file of user study"
#
# Changes applied to Sheet1:
#  - B3 value  " Mike Collins"   -> "M. C."   (anonymize the one real name
#    still present in the data)
#  - B1 header "Name"            -> "ID"      (column no longer holds a name)
#  - Column B narrower now that the long name text is gone
#  - Active selection moved from the old scratch selection (M1:M1048576)
#    to B2, matching the frozen-pane view of the data area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order matters for shared-string table layout: anonymize the data cell
# before renaming the header so "M. C." is registered ahead of "ID".
$ws.Range("B3").Value = "M. C."
$ws.Range("B1").Value = "ID"

# Narrow column B now that the long name text is gone (closest reachable
# width to the authored 14.42578125).
$ws.Columns.Item(2).ColumnWidth = 13.7109375

# Update the stored selection for the frozen/bottom-right pane to B2.
$ws.Range("B2").Select()
